$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 301, shifting existing rows 301-402 down to 302-403
$ws.Rows("301:301").Insert()

# Populate the new row 301 with the new data record
$ws.Range("A301").Value = 3
$ws.Range("B301").Value = "Femacal de La Calera"
$ws.Range("C301").Value = "Coquimbo"
$ws.Range("D301").Value = 44809
$ws.Range("D301").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E301").Value = 5
$ws.Range("F301").Value = 100112009
$ws.Range("G301").Value = "Acelga"
$ws.Range("H301").Value = "Sin especificar"
$ws.Range("I301").Value = "Primera"
$ws.Range("J301").Value = 230
$ws.Range("K301").Value = 3300
$ws.Range("L301").Value = 3500
$ws.Range("M301").Value = 3404
$ws.Range("N301").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O301").Value = "Provincia de Quillota"
$ws.Range("P301").Value = 567
$ws.Range("Q301").Value = 6
$ws.Range("R301").Value = "Hortaliza"
